$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K header: "idhorario" (added to sharedStrings automatically),
# right-aligned like the other header cells (style index 1).
$ws.Cells.Item(1, 11).Value = "idhorario"
$ws.Cells.Item(1, 11).HorizontalAlignment = -4152

# Populate column K (idhorario) for every data row, grouped in the same
# bands used by the source data (1 for rows 2-61, 2 for rows 62-121, etc.)
for ($i = 2; $i -le 61; $i++)   { $ws.Cells.Item($i, 11).Value = 1 }
for ($i = 62; $i -le 121; $i++) { $ws.Cells.Item($i, 11).Value = 2 }
for ($i = 122; $i -le 161; $i++) { $ws.Cells.Item($i, 11).Value = 3 }
for ($i = 162; $i -le 221; $i++) { $ws.Cells.Item($i, 11).Value = 4 }
for ($i = 222; $i -le 255; $i++) { $ws.Cells.Item($i, 11).Value = 5 }

# Match the new selection left behind in the saved workbook.
$ws.Range("G64").Select() | Out-Null
